$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price/Volume columns so values like "35.116.61"
# or "0.0693" are stored as literal text rather than being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.116.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.31"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.20%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.18"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.17%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.122.36"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.856.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.086.89"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.96"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.68"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +23.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.124"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0556"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.99"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +27.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.840"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +20.05%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.08"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.31%  "
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "90.36"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0201"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.340.18"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.99"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.07%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +43.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0556"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.35%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.60"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.037.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0681"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.37%  "
